# Fix mistake regarding H2 transport:
# On the "Fueltrade" sheet, column F (Comm2) was incorrectly set to "H2GC"
# for rows 4-11; it should be "H2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

for ($r = 4; $r -le 11; $r++) {
    $ws.Cells.Item($r, 6).Value = "H2"
}

# Update the active selection to match the recorded view state.
$ws.Range("H14").Select()
